$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.563.04'
$ws.Range("E2").Value = '  +1.81%  '

$ws.Range("D3").Value = '1.664.64'
$ws.Range("E3").Value = '  +0.88%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.13'
$ws.Range("E5").Value = '  -0.56%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4797'
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2614'
$ws.Range("E8").Value = '  -0.63%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06145'
$ws.Range("E9").Value = '  +1.82%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07090'
$ws.Range("E10").Value = '  -0.06%  '

$ws.Range("D11").Value = '1.661.76'
$ws.Range("E11").Value = '  +0.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.75'
$ws.Range("E12").Value = '  +1.86%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5890'
$ws.Range("E13").Value = '  -4.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.370'
$ws.Range("E14").Value = '  -4.38%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '74.47'
$ws.Range("E15").Value = '  +1.90%  '

$ws.Range("E16").Value = '  +0.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9997'
$ws.Range("E17").Value = '  -0.03%  '

$ws.Range("D18").Value = '25.547.87'
$ws.Range("E18").Value = '  +1.80%  '

$ws.Range("E19").Value = '  +2.84%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.39'
$ws.Range("E20").Value = '  +0.16%  '

$ws.Range("D21").Value = '1.882.89'
$ws.Range("E21").Value = '  +1.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.426'
$ws.Range("E22").Value = '  +0.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.658'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.302'
$ws.Range("E24").Value = '  +1.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.35'
$ws.Range("E25").Value = '  +0.35%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.04'
$ws.Range("E26").Value = '  +2.05%  '

$ws.Range("E27").Value = '  +0.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '104.66'
$ws.Range("E28").Value = '  +2.85%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.686'
$ws.Range("E29").Value = '  +0.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.955'
$ws.Range("E30").Value = '  +4.25%  '

$ws.Range("E31").Value = '  +2.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07628'
$ws.Range("E32").Value = '  -3.70%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9992'
$ws.Range("E33").Value = '  +0.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04283'
$ws.Range("E34").Value = '  -5.58%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.618'
$ws.Range("E35").Value = '  +0.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6127'
$ws.Range("E36").Value = '  +5.85%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9484'
$ws.Range("E37").Value = '  +0.77%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.610'
$ws.Range("E38").Value = '  -0.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8461'
$ws.Range("E39").Value = '  +0.90%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9998'
$ws.Range("E40").Value = '  +0.04%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01496'
$ws.Range("E41").Value = '  -2.74%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.869'
$ws.Range("E42").Value = '  +3.18%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.68'
$ws.Range("E43").Value = '  -1.07%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3758'
$ws.Range("E44").Value = '  +1.55%  '

$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.700'
$ws.Range("E45").Value = '  -1.97%  '

$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1119'
$ws.Range("E46").Value = '  -0.88%  '

$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.213'
$ws.Range("E47").Value = '  +3.08%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05259'
$ws.Range("E48").Value = '  +1.74%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.37'
$ws.Range("E49").Value = '  -0.82%  '

$ws.Range("B50").Value = 'TrueUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.001'
$ws.Range("E50").Value = '  +0.07%  '

$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("E51").Value = '  +0.27%  '
